$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '69.781.81'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  -0.92%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.491.95'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  -1.83%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '606.13'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +3.14%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '192.44'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +2.16%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.625'
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  +0.84%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.08%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.212'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -1.55%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.659'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +2.16%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '53.37'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -1.49%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.0000306'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -1.37%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '9.59'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +1.96%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '4.059.09'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -1.62%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '614.74'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +5.57%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '69.918.32'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -0.79%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '12.66'
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -0.77%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '18.82'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -0.85%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '3.492.67'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -1.43%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -0.16%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '0.989'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -0.49%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '17.76'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -0.98%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '105.07'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +11.10%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '4.64'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -0.11%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '5.05'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +3.98%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +4.07%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '10.95'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -0.50%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '9.88'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +5.57%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '34.23'
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +5.91%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '7.14'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +1.26%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '4.31'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +13.33%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '12.62'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +3.53%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.114'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -0.25%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '64.22'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +1.29%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '3.741.70'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +1.56%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'Dai'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.999'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -0.07%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'Bittensor'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '520.25'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -1.32%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'Fetch.AI'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '3.06'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -5.01%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.0₃0792'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +0.58%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.390'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -3.96%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '3.57'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +1.00%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '36.59'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -3.53%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -0.43%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +1.60%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -2.73%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +2.01%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -3.97%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.33%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '8.74'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -5.30%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '132.74'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -2.05%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -5.95%  '
